# Hemos cambiado la formula de Ventas objetivo
# Update "uds. Objetivo semana pasada" (R), "Tendencia Consumo" (T)
# and "Diferencia Stock" (L) values across the order sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R (uds. Objetivo semana pasada)
$ws.Range("R5").Value = 3
$ws.Range("R6").Value = 1
$ws.Range("R8").Value = 2
$ws.Range("R9").Value = 1
$ws.Range("R10").Value = 4
$ws.Range("R11").Value = 1
$ws.Range("R13").Value = 2
$ws.Range("R16").Value = 1
$ws.Range("R22").Value = 2
$ws.Range("R23").Value = 2
$ws.Range("R25").Value = 4
$ws.Range("R37").Value = 6
$ws.Range("R39").Value = 2
$ws.Range("R47").Value = 1

# Column T (Tendencia Consumo)
$ws.Range("T9").Value = 0

# Column L (Diferencia Stock)
$ws.Range("L10").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("L37").Value = 0
